# Getting started with PHP
# Fill in the Topic (column C) for rows 87-90 and adjust row heights for rows 102-106

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C87").Value = "Fluid Layout"
$ws.Range("C88").Value = "Media Query"
$ws.Range("C89").Value = "Media Query"
$ws.Range("C90").Value = "Getting Started with PHP"

$ws.Rows.Item(102).RowHeight = 19.5
$ws.Rows.Item(103).RowHeight = 19.5
$ws.Rows.Item(104).RowHeight = 19.5
$ws.Rows.Item(105).RowHeight = 19.5
$ws.Rows.Item(106).RowHeight = 20.25
